$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.240179
$ws.Range("H2").Value = 3.720537
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.695610666666666
$ws.Range("N2").Value = 14.086832
$ws.Range("O2").Value = 0.1802066564018305
$ws.Range("P2").Value = 0.1802066564018305
$ws.Range("Q2").Value = 5.823397740975999
$ws.Range("R2").Value = 52.410579668784
$ws.Range("S2").Value = 0.1802066564018305
$ws.Range("T2").Value = 0.1802066564018305

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.240179
$ws.Range("H3").Value = 3.720537
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 15.51448033333333
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5954098039960916
$ws.Range("P3").Value = 0.5954098039960916
$ws.Range("Q3").Value = 19.240732705313
$ws.Range("R3").Value = 173.166594347817
$ws.Range("S3").Value = 0.5954098039960916
$ws.Range("T3").Value = 0.5954098039960916

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.240179
$ws.Range("H4").Value = 3.720537
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.846719333333333
$ws.Range("N4").Value = 17.540158
$ws.Range("O4").Value = 0.2243835396020779
$ws.Range("P4").Value = 0.2243835396020779
$ws.Range("Q4").Value = 7.250978536093998
$ws.Range("R4").Value = 65.25880682484599
$ws.Range("S4").Value = 0.2243835396020779
$ws.Range("T4").Value = 0.2243835396020779
